# The workbook gained one new weekly price record for "Perejil" at
# Feria Lagunitas de Puerto Montt. It was inserted as a new row 377,
# which pushes every following record down by one row (old row 377
# becomes 378, ..., old row 452 becomes the new row 453). The sheet's
# used range therefore grows from A1:R452 to A1:R453.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 377; Excel shifts rows 377:452 down to 378:453.
$ws.Rows.Item(377).Insert()

# Populate the newly inserted row 377 with the new record's data.
$ws.Cells.Item(377, 1).Value2  = 4
$ws.Cells.Item(377, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(377, 3).Value2  = "Los Lagos"
$ws.Cells.Item(377, 4).Value2  = 45204
$ws.Cells.Item(377, 5).Value2  = 10
$ws.Cells.Item(377, 6).Value2  = 100112044
$ws.Cells.Item(377, 7).Value2  = "Perejil"
$ws.Cells.Item(377, 8).Value2  = "Sin especificar"
$ws.Cells.Item(377, 9).Value2  = "Primera"
$ws.Cells.Item(377, 10).Value2 = 80
$ws.Cells.Item(377, 11).Value2 = 6000
$ws.Cells.Item(377, 12).Value2 = 6000
$ws.Cells.Item(377, 13).Value2 = 6000
$ws.Cells.Item(377, 14).Value2 = "$/docena de atados (3 kilos)"
$ws.Cells.Item(377, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(377, 16).Value2 = 2000
$ws.Cells.Item(377, 17).Value2 = 3
$ws.Cells.Item(377, 18).Value2 = "Hortaliza"
